$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply cell updates as plain text (matches source workbook which stores
# all data-row cells as inline/shared strings, even numeric-looking values).
# Trick: force the NumberFormat to Text ("@") before the assignment so the
# value is not auto-coerced to a Number/Percentage, then reset the cell
# style back to "Normal" so no stray style index is left behind.
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "254.97"
Set-TextValue "E2" "3.58%"
Set-TextValue "G2" "10"
Set-TextValue "D3" "28.24"
Set-TextValue "E3" "-6.76%"
Set-TextValue "G3" "10"
Set-TextValue "D4" "5.254"
Set-TextValue "E4" "1.99%"
Set-TextValue "G4" "10"
Set-TextValue "D5" "0.05855"
Set-TextValue "E5" "1.40%"
Set-TextValue "G5" "10"
Set-TextValue "D6" "6.727"
Set-TextValue "E6" "1.02%"
Set-TextValue "G6" "10"
Set-TextValue "D7" "0.8662"
Set-TextValue "E7" "2.00%"
Set-TextValue "G7" "10"
Set-TextValue "D8" "1.004"
Set-TextValue "E8" "17.21%"
Set-TextValue "G8" "10"
Set-TextValue "D9" "0.1408"
Set-TextValue "E9" "1.22%"
Set-TextValue "G9" "10"
Set-TextValue "D10" "0.07168"
Set-TextValue "G10" "10"
Set-TextValue "D11" "0.03171"
Set-TextValue "E11" "-2.69%"
Set-TextValue "G11" "10"
Set-TextValue "E12" "-1.52%"
Set-TextValue "G12" "10"
Set-TextValue "E13" "1.27%"
Set-TextValue "G13" "10"
Set-TextValue "D14" "0.0006072"
Set-TextValue "E14" "-94.05%"
Set-TextValue "G14" "10"
Set-TextValue "D15" "0.005808"
Set-TextValue "E15" "-4.06%"
Set-TextValue "G15" "10"
Set-TextValue "D16" "3.500"
Set-TextValue "E16" "-0.68%"
Set-TextValue "G16" "10"
Set-TextValue "D17" "3.225"
Set-TextValue "E17" "-0.59%"
Set-TextValue "G17" "10"
Set-TextValue "D18" "2.203"
Set-TextValue "E18" "0.75%"
Set-TextValue "G18" "10"
Set-TextValue "D19" "0.3178"
Set-TextValue "E19" "0.43%"
Set-TextValue "G19" "10"
Set-TextValue "D20" "0.03472"
Set-TextValue "E20" "2.51%"
Set-TextValue "G20" "10"
Set-TextValue "E21" "-0.61%"
Set-TextValue "G21" "10"
Set-TextValue "D22" "3.534"
Set-TextValue "E22" "1.34%"
Set-TextValue "G22" "10"
Set-TextValue "D23" "0.04151"
Set-TextValue "E23" "0.45%"
Set-TextValue "G23" "10"
Set-TextValue "D24" "0.1381"
Set-TextValue "E24" "-2.00%"
Set-TextValue "G24" "10"
Set-TextValue "D25" "0.001226"
Set-TextValue "E25" "0.06%"
Set-TextValue "G25" "10"
Set-TextValue "D26" "0.004806"
Set-TextValue "E26" "15.77%"
Set-TextValue "G26" "10"
Set-TextValue "E27" "0.08%"
Set-TextValue "G27" "10"
Set-TextValue "E28" "1.23%"
Set-TextValue "G28" "10"
Set-TextValue "G29" "10"
Set-TextValue "G30" "10"
Set-TextValue "G31" "10"
Set-TextValue "G32" "10"
Set-TextValue "G33" "10"
Set-TextValue "G34" "10"
Set-TextValue "G35" "10"
Set-TextValue "G36" "10"
Set-TextValue "G37" "10"
Set-TextValue "G38" "10"
Set-TextValue "G39" "10"
Set-TextValue "D40" "0.03804"
Set-TextValue "E40" "1.16%"
Set-TextValue "G40" "10"
Set-TextValue "B41" "BKEXToken"
Set-TextValue "C41" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1101"
Set-TextValue "E41" "2.62%"
Set-TextValue "G41" "10"
Set-TextValue "B42" "KickToken"
Set-TextValue "C42" "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D42" "0.003813"
Set-TextValue "E42" "-32.62%"
Set-TextValue "G42" "10"
Set-TextValue "D43" "0.002339"
Set-TextValue "E43" "-5.29%"
Set-TextValue "G43" "10"
Set-TextValue "D44" "0.009696"
Set-TextValue "E44" "-2.59%"
Set-TextValue "G44" "10"
Set-TextValue "D45" "0.00005237"
Set-TextValue "E45" "-4.42%"
Set-TextValue "G45" "10"
Set-TextValue "D46" "0.00000000750"
Set-TextValue "E46" "0.05%"
Set-TextValue "G46" "10"
Set-TextValue "D47" "0.09302"
Set-TextValue "E47" "31.06%"
Set-TextValue "G47" "10"
Set-TextValue "D48" "0.002151"
Set-TextValue "E48" "-12.78%"
Set-TextValue "G48" "10"
Set-TextValue "D49" "0.00002101"
Set-TextValue "E49" "0.05%"
Set-TextValue "G49" "10"
Set-TextValue "D50" "0.0002000"
Set-TextValue "E50" "0.05%"
Set-TextValue "G50" "10"
Set-TextValue "G51" "10"
